$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric are stored as text, matching the source data
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.198.09'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.848.88'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").Value = '0.9983'
$ws.Range("D5").Value = '245.79'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").Value = '0.6983'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = '0.9991'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.07715'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '0.3064'
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").Value = '23.60'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").Value = '0.07819'
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '93.24'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '1.841.80'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").Value = '5.129'
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D15").Value = '0.6861'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '6.658'
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("D17").Value = '0.000008308'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '29.178.49'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '241.93'
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("D20").Value = '2.080.40'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").Value = '12.77'
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("D22").Value = '0.9989'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '7.506'
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '0.1519'
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("D26").Value = '159.60'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '8.834'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").Value = '18.27'
$ws.Range("E28").Value = '  -1.56%  '
$ws.Range("D29").Value = '1.543'
$ws.Range("E29").Value = '  -1.43%  '
$ws.Range("D30").Value = '4.243'
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").Value = '4.189'
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = '0.05117'
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("D34").Value = '0.7853'
$ws.Range("E34").Value = '  +3.73%  '
$ws.Range("D35").Value = '1.865'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = '1.148'
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").Value = '2.692'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = '1.323.59'
$ws.Range("E38").Value = '  +8.25%  '
$ws.Range("D39").Value = '0.01870'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").Value = '0.9479'
$ws.Range("E41").Value = '  +5.13%  '
$ws.Range("D42").Value = '6.056'
$ws.Range("E42").Value = '  +5.12%  '
$ws.Range("D43").Value = '107.78'
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").Value = '0.9989'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").Value = '9.725'
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("D47").Value = '0.5177'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = '1.982.74'
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").Value = '64.38'
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").Value = '1.764'
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '6.995'
$ws.Range("E51").Value = '  -0.47%  '
